$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AltaCuentaCTS")

# Rename the "cuenta" header in H1 to "Tipo de Producto"
$ws.Range("H1").Value = "Tipo de Producto"

# Move the active selection to J9
$ws.Range("J9").Select()
